$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "before" values in column A and "after" values in column B.
# Formulas in columns C and E recalculate automatically.
$ws.Range("A1").Value = 866
$ws.Range("B1").Value = 1670

$ws.Range("B2").Value = 1150

$ws.Range("A4").Value = 1743
$ws.Range("B4").Value = 2328

$ws.Range("A5").Value = 1150
$ws.Range("B5").Value = 1400

$ws.Range("A7").Value = 2401
$ws.Range("B7").Value = 2839

$ws.Range("A8").Value = 1400
$ws.Range("B8").Value = 1800

$ws.Range("A10").Value = 2912
$ws.Range("B10").Value = 3351

$ws.Range("A11").Value = 1800
$ws.Range("B11").Value = 2300

$ws.Range("A13").Value = 3424

$ws.Range("A14").Value = 2300

# Update the active selection to match the saved view state.
$ws.Range("E8").Select()
